# Update cryptos list: refresh Price (column D) and Volume(1h) (column E) values.
# NumberFormat is forced to Text first so Excel stores these as literal strings
# instead of re-interpreting numeric-looking values (dropping trailing zeros, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.700.18"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.723.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9986"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4820"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2577"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06167"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.720.11"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.84"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6026"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.450"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.00"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9988"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.529.36"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9987"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007136"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.941.73"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.404"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.540"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.037"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.57"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.24"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.768"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.73"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.010"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07898"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04489"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.595"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9976"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6157"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9310"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.000"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.444"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.37%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.599"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.77"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3814"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.752"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05359"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.893"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.96"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.238"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.17"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.70%  "
